$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.727.95"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "2.641.73"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.106"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.12%  "
$ws.Range("D15").Value = "60.703.73"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.528"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +9.01%  "
$ws.Range("D29").Value = "0.0₃0804"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +6.93%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("E36").Value = "  +6.75%  "
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "339.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.43%  "
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.907"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.68%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0249"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0562"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "2.086.88"
$ws.Range("E51").Value = "  +1.81%  "
